$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.757.62"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "2.337.96"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("D5").Value = "'239.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("E6").Value = "  -3.92%  "
$ws.Range("D7").Value = "'72.08"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.16%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -6.46%  "
$ws.Range("D10").Value = "'0.0990"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.37%  "
$ws.Range("D11").Value = "'58.10"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").Value = "'32.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.52%  "
$ws.Range("D13").Value = "'0.108"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").Value = "'7.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.72%  "
$ws.Range("D15").Value = "2.685.65"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "'16.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.01%  "
$ws.Range("D17").Value = "'0.897"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.89%  "
$ws.Range("D18").Value = "2.337.90"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").Value = "43.676.85"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("E20").Value = "  -2.79%  "
$ws.Range("D21").Value = "'77.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "'6.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").Value = "'251.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "'1.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.38%  "
$ws.Range("E26").Value = "  +2.52%  "
$ws.Range("D27").Value = "'2.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").Value = "'10.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.84%  "
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("D30").Value = "'176.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("E31").Value = "  -4.17%  "
$ws.Range("E32").Value = "  -2.12%  "
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("D34").Value = "'0.0735"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("D35").Value = "'5.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("E37").Value = "  -2.86%  "
$ws.Range("E38").Value = "  -1.77%  "
$ws.Range("E39").Value = "  -3.16%  "
$ws.Range("D40").Value = "'5.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +26.19%  "
$ws.Range("E41").Value = "  -3.10%  "
$ws.Range("D42").Value = "'65.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +19.29%  "
$ws.Range("D43").Value = "'9.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.74%  "
$ws.Range("E44").Value = "  +4.06%  "
$ws.Range("D45").Value = "'18.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.96%  "
$ws.Range("E46").Value = "  -3.84%  "
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").Value = "'1.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.00%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'2.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.69%  "
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "'2.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.40%  "
$ws.Range("D51").Value = "'97.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.58%  "
